$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per row and participate in the cyclic shift
$cols = @("A","B","D","E","F","G","H","I","Q","R","S","Z","AB")

# Capture current values for rows 2-5 for each relevant column
$values = @{}
foreach ($col in $cols) {
    $values[$col] = @{}
    for ($r = 2; $r -le 5; $r++) {
        $addr = "$col$r"
        $values[$col][$r] = $ws.Range($addr).Value()
    }
}

# Write back with cyclic shift: new row2 = old row3, new row3 = old row4,
# new row4 = old row5, new row5 = old row2
$srcMap = @{2=3; 3=4; 4=5; 5=2}
foreach ($col in $cols) {
    foreach ($destRow in 2..5) {
        $srcRow = $srcMap[$destRow]
        $val = $values[$col][$srcRow]
        $ws.Range("$col$destRow").Value = $val
    }
}
